$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# ALC row 16
$ws.Range("H16").Value = 8000
$ws.Range("J16").Value = 8000
$ws.Range("L16").Value = 8000
$ws.Range("N16").Value = -8460

# ALC row 32
$ws.Range("H32").Value = 3067
$ws.Range("J32").Value = 2725.5
$ws.Range("L32").Value = 2725.5
$ws.Range("N32").Value = -3377.5

# ALC row 62
$ws.Range("H62").Value = 2382.9644
$ws.Range("I62").Value = 2218.5625
$ws.Range("J62").Value = 2602.1667
$ws.Range("K62").Value = 2218.5625
$ws.Range("L62").Value = 2602.1667
$ws.Range("M62").Value = -1594.5625
$ws.Range("N62").Value = -3850.1667

# ALC row 65
$ws.Range("H65").Value = 2382.9644
$ws.Range("I65").Value = 2218.5625
$ws.Range("J65").Value = 2602.1667
$ws.Range("K65").Value = 11092.8125
$ws.Range("L65").Value = 13010.8335
$ws.Range("M65").Value = -7972.8125
$ws.Range("N65").Value = -19250.8335

$ws = $wb.Worksheets.Item("ARM")
# ARM row 2
$ws.Range("H2").Value = 1410.25
$ws.Range("I2").Value = 1510
$ws.Range("J2").Value = 911.5
$ws.Range("K2").Value = 1510
$ws.Range("L2").Value = 911.5
$ws.Range("M2").Value = -1397
$ws.Range("N2").Value = -1137.5

# ARM row 9
$ws.Range("H9").Value = 19950
$ws.Range("J9").Value = 19950
$ws.Range("L9").Value = 19950
$ws.Range("N9").Value = -20290

# ARM row 20
$ws.Range("H20").Value = 19950
$ws.Range("J20").Value = 19950
$ws.Range("L20").Value = 19950
$ws.Range("N20").Value = -20490

# ARM row 32
$ws.Range("H32").Value = 14711995
$ws.Range("I32").Value = 21741378
$ws.Range("J32").Value = 14196
$ws.Range("K32").Value = 21741378
$ws.Range("L32").Value = 14196
$ws.Range("M32").Value = -21741091
$ws.Range("N32").Value = -14770

# ARM row 45
$ws.Range("H45").Value = 1933.2759
$ws.Range("I45").Value = 1725.48
$ws.Range("J45").Value = 3232
$ws.Range("K45").Value = 1725.48
$ws.Range("L45").Value = 3232
$ws.Range("M45").Value = -1348.48
$ws.Range("N45").Value = -3986

# ARM row 74
$ws.Range("H74").Value = 3457.5854
$ws.Range("I74").Value = 4449.4644
$ws.Range("J74").Value = 1321.2307
$ws.Range("K74").Value = 4449.4644
$ws.Range("L74").Value = 1321.2307
$ws.Range("M74").Value = -3575.4644
$ws.Range("N74").Value = -3069.2307

# ARM row 77
$ws.Range("H77").Value = 3457.5854
$ws.Range("I77").Value = 4449.4644
$ws.Range("J77").Value = 1321.2307
$ws.Range("K77").Value = 22247.322
$ws.Range("L77").Value = 6606.1535
$ws.Range("M77").Value = -17879.322
$ws.Range("N77").Value = -15342.1535

# ARM row 86
$ws.Range("H86").Value = 37666.668
$ws.Range("I86").Value = 12000
$ws.Range("K86").Value = 12000
$ws.Range("M86").Value = -10814

# ARM row 89
$ws.Range("H89").Value = 37666.668
$ws.Range("I89").Value = 12000
$ws.Range("K89").Value = 36000
$ws.Range("M89").Value = -30072

# ARM row 116
$ws.Range("H116").Value = 1410.25
$ws.Range("I116").Value = 1510
$ws.Range("J116").Value = 911.5
$ws.Range("K116").Value = 1510
$ws.Range("L116").Value = 911.5
$ws.Range("M116").Value = 784
$ws.Range("N116").Value = -5499.5

# ARM row 132
$ws.Range("H132").Value = 3149.1072
$ws.Range("I132").Value = 2401.8948
$ws.Range("J132").Value = 4726.5557
$ws.Range("K132").Value = 7205.6844
$ws.Range("L132").Value = 14179.6671
$ws.Range("M132").Value = -4675.6844
$ws.Range("N132").Value = -19239.6671

$ws = $wb.Worksheets.Item("BSM")
# BSM row 3
$ws.Range("H3").Value = 1410.25
$ws.Range("I3").Value = 1510
$ws.Range("J3").Value = 911.5
$ws.Range("K3").Value = 1510
$ws.Range("L3").Value = 911.5
$ws.Range("M3").Value = -1396
$ws.Range("N3").Value = -1139.5

$ws = $wb.Worksheets.Item("CRP")
# CRP row 107
$ws.Range("H107").Value = 435592.53
$ws.Range("I107").Value = 909790.8
$ws.Range("J107").Value = 910.75
$ws.Range("K107").Value = 909790.8
$ws.Range("L107").Value = 910.75
$ws.Range("M107").Value = -907870.8
$ws.Range("N107").Value = -4750.75

# CRP row 122
$ws.Range("H122").Value = 1304.6154
$ws.Range("I122").Value = 1182.5
$ws.Range("J122").Value = 1500
$ws.Range("K122").Value = 3547.5
$ws.Range("L122").Value = 4500
$ws.Range("M122").Value = -1097.5
$ws.Range("N122").Value = -9400

# CRP row 132
$ws.Range("H132").Value = 125003240
$ws.Range("I132").Value = 200001810
$ws.Range("J132").Value = 5637.6665
$ws.Range("K132").Value = 600005430
$ws.Range("L132").Value = 16912.9995
$ws.Range("M132").Value = -600002900
$ws.Range("N132").Value = -21972.9995

$ws = $wb.Worksheets.Item("CUL")
# CUL row 5
$ws.Range("H5").Value = 583.04
$ws.Range("I5").Value = 554.4545000000001
$ws.Range("K5").Value = 1663.3635
$ws.Range("M5").Value = -1551.3635

# CUL row 13
$ws.Range("H13").Value = 2396.6667
$ws.Range("I13").Value = 1480
$ws.Range("J13").Value = 2580
$ws.Range("K13").Value = 4440
$ws.Range("L13").Value = 7740
$ws.Range("M13").Value = -4272
$ws.Range("N13").Value = -8076

# CUL row 135
$ws.Range("H135").Value = 583.04
$ws.Range("I135").Value = 554.4545000000001
$ws.Range("K135").Value = 4990.0905
$ws.Range("M135").Value = -2455.0905

$ws = $wb.Worksheets.Item("GSM")
# GSM row 122
$ws.Range("H122").Value = 2986.524
$ws.Range("I122").Value = 2754.3845
$ws.Range("J122").Value = 3363.75
$ws.Range("K122").Value = 8263.1535
$ws.Range("L122").Value = 10091.25
$ws.Range("M122").Value = -5813.1535
$ws.Range("N122").Value = -14991.25

# GSM row 132
$ws.Range("H132").Value = 5615.077
$ws.Range("I132").Value = 5749.2812
$ws.Range("J132").Value = 5001.5713
$ws.Range("K132").Value = 17247.8436
$ws.Range("L132").Value = 15004.7139
$ws.Range("M132").Value = -14717.8436
$ws.Range("N132").Value = -20064.7139

$ws = $wb.Worksheets.Item("LTW")
# LTW row 40
$ws.Range("H40").Value = 4658.8887
$ws.Range("I40").Value = 3316.6667
$ws.Range("J40").Value = 5330
$ws.Range("K40").Value = 3316.6667
$ws.Range("L40").Value = 5330
$ws.Range("M40").Value = -3180.6667
$ws.Range("N40").Value = -5602

# LTW row 136
$ws.Range("H136").Value = 1491.5
$ws.Range("I136").Value = 970.3125
$ws.Range("J136").Value = 3576.25
$ws.Range("K136").Value = 2910.9375
$ws.Range("L136").Value = 10728.75
$ws.Range("M136").Value = -360.9375
$ws.Range("N136").Value = -15828.75

$ws = $wb.Worksheets.Item("WVR")
# WVR row 107
$ws.Range("H107").Value = 737.25
$ws.Range("I107").Value = 760
$ws.Range("K107").Value = 2280
$ws.Range("M107").Value = -360

# WVR row 132
$ws.Range("H132").Value = 1961.6538
$ws.Range("I132").Value = 1633.3334
$ws.Range("J132").Value = 2700.375
$ws.Range("K132").Value = 4900.0002
$ws.Range("L132").Value = 8101.125
$ws.Range("M132").Value = -2370.0002
$ws.Range("N132").Value = -13161.125
